$d = $word.ActiveDocument

# Helper: return the 1-based index of the paragraph that contains character
# position $pos.
function Get-ParagraphIndexAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return $doc.Paragraphs.Count
}

# ---------------------------------------------------------------------------
# 1. Merge the split "Requirements - AudioAttributes" title runs into a
#    single run (this also removes the w:proofErr spell-check markers that
#    Word had placed around "AudioAttributes").
# ---------------------------------------------------------------------------
$dash = [char]8211
$titleText = "Requirements " + $dash + " AudioAttributes"

$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("AudioAttributes", $false, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $titleIndex = Get-ParagraphIndexAt $d $searchRange.Start
} else {
    $titleIndex = 1
}
$titleRange = $d.Paragraphs.Item($titleIndex).Range

$titleXml = '<?xml version="1.0"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p w:rsidR="00690A57" w:rsidRDefault="002F7D97" w:rsidP="002F7D97">' +
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2. After the "Write Audio" bullet, add three new sub-bullets describing
#    the channel-merging requirement.
# ---------------------------------------------------------------------------
$searchRange2 = $d.Content.Duplicate
$found2 = $searchRange2.Find.Execute("Write Audio", $true, $false, $false,
                                      $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $waIndex = Get-ParagraphIndexAt $d $searchRange2.Start
} else {
    $waIndex = 5
}
$waRange = $d.Paragraphs.Item($waIndex).Range

$newItemsXml = '<?xml version="1.0"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p w:rsidR="00C74F5E" w:rsidRDefault="00C74F5E" w:rsidP="00C6427E">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Write Audio</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Merge channels:</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Set levels for channels</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Retrofit normalisation factor from channel into other channels of audio input.</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$waRange.InsertXML($newItemsXml)
